# Apply weekly data refresh to "Fruta, Comercializadora del Agro de Limarí - Damasco"
# The underlying market rows (2-7) were updated/reordered to reflect the new week's
# price report. Only the cells that actually change values are touched below.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was Dina / Primera / Coquimbo-Metropolitana, now Dina / Especial / O'Higgins)
$ws.Range("D2").Value = 44189
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 23500
$ws.Range("O2").Value = 24000
$ws.Range("P2").Value = 23750
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 1319
$ws.Range("T2").Value = 18

# Row 3 (was Dina / Segunda, now Dina / Primera / O'Higgins)
$ws.Range("D3").Value = 44189
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 21500
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21750
$ws.Range("Q3").Value = "$/caja 18 kilos"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 1208
$ws.Range("T3").Value = 18

# Row 4 (was Castle Brite, now Dina)
$ws.Range("D4").Value = 44161
$ws.Range("K4").Value = "Dina"
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20500
$ws.Range("P4").Value = 20250
$ws.Range("Q4").Value = "$/caja 15 kilos"
$ws.Range("S4").Value = 1350
$ws.Range("T4").Value = 15

# Row 5 (was Castle Brite / Primera, now Dina / Segunda)
$ws.Range("D5").Value = 44161
$ws.Range("K5").Value = "Dina"
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 18000
$ws.Range("O5").Value = 18500
$ws.Range("P5").Value = 18250
$ws.Range("S5").Value = 1217

# Row 6 (was Dina / Especial / O'Higgins, now Castle Brite / Primera / Metropolitana)
$ws.Range("D6").Value = 44160
$ws.Range("K6").Value = "Castle Brite"
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 240
$ws.Range("N6").Value = 20500
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20750
$ws.Range("Q6").Value = "$/caja 15 kilos"
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 1383
$ws.Range("T6").Value = 15

# Row 7 (was Dina / Primera / O'Higgins, now Castle Brite / Primera / Metropolitana)
$ws.Range("D7").Value = 44175
$ws.Range("K7").Value = "Castle Brite"
$ws.Range("M7").Value = 300
$ws.Range("N7").Value = 21000
$ws.Range("P7").Value = 21500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 1194
